# DN 4.0 - Java FSE Mandatory hands-on detail.xlsx
# "adding week 6 react HOL exercises and images"
#
# Visible, reproducible effects of this commit on the worksheet data:
#   - Rows 31:37 in column H flip their status from "Pending" to "Completed"
#     (the week-6 React HOL rows now have hands-on material attached).
#   - The view has scrolled down (frozen header still covers rows 1-2) and
#     the active/selected cell moved from I30 to I37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Mark the week 6 React HOL rows (31-37) as "Completed" in the Status column.
$ws.Range("H31:H37").Value = "Completed"

# Move the selection/view down to where the newly completed rows are.
$ws.Range("I37").Select()
